$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the role for this user (E2): role id 1 -> 2
$ws.Range("E2").Value = 2

# Update the active cell/selection to E3 (as left after the edit)
$ws.Range("E3").Select()
